$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Jeff"
$ws.Range("A3").Value = "Sturek"
$ws.Range("A4").Value = "doin stuff"

$ws.Range("F7").Select()
